# This script updates the "想去人数" (want-to-go count) values in column F
# across the "展览", "演出" and "全部类型" worksheets, matching the data
# refresh captured in the commit "Update gh-pages to output generated at
# 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (sheetId=1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 5190
$ws1.Range("F5").Value  = 7494
$ws1.Range("F7").Value  = 72
$ws1.Range("F9").Value  = 603
$ws1.Range("F12").Value = 4336
$ws1.Range("F14").Value = 108
$ws1.Range("F15").Value = 110
$ws1.Range("F16").Value = 2935
$ws1.Range("F19").Value = 213
$ws1.Range("F20").Value = 513
$ws1.Range("F22").Value = 465
$ws1.Range("F23").Value = 316
$ws1.Range("F26").Value = 1196
$ws1.Range("F27").Value = 95
$ws1.Range("F28").Value = 1391
$ws1.Range("F30").Value = 585
$ws1.Range("F35").Value = 107
$ws1.Range("F37").Value = 2938
$ws1.Range("F40").Value = 93
$ws1.Range("F42").Value = 50

# --- Sheet: 演出 (sheetId=2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 14

# --- Sheet: 全部类型 (sheetId=4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 5190
$ws4.Range("F5").Value  = 7494
$ws4.Range("F7").Value  = 72
$ws4.Range("F9").Value  = 603
$ws4.Range("F12").Value = 4336
$ws4.Range("F14").Value = 108
$ws4.Range("F15").Value = 110
$ws4.Range("F16").Value = 2935
$ws4.Range("F19").Value = 213
$ws4.Range("F20").Value = 513
$ws4.Range("F22").Value = 465
$ws4.Range("F24").Value = 316
$ws4.Range("F27").Value = 1196
$ws4.Range("F28").Value = 95
$ws4.Range("F29").Value = 1391
$ws4.Range("F31").Value = 585
$ws4.Range("F36").Value = 107
$ws4.Range("F38").Value = 2938
$ws4.Range("F39").Value = 14
$ws4.Range("F42").Value = 93
$ws4.Range("F44").Value = 50

$wb.Save()
